$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 06:36"

# Update India (row 5)
$ws.Range("B5").Value = 6394068
$ws.Range("C5").Value = 2108
$ws.Range("D5").Value = 5352078
$ws.Range("E5").Value = 942186

# Update Pakistan (row 24)
$ws.Range("B24").Value = 313431
$ws.Range("C24").Value = 625
$ws.Range("D24").Value = 298055
$ws.Range("E24").Value = 8877
$ws.Range("G24").Value = 15
$ws.Range("H24").Value = 6499

# Update Tailandia (row 141)
$ws.Range("B141").Value = 3575
$ws.Range("C141").Value = 6
$ws.Range("D141").Value = 3384
$ws.Range("E141").Value = 132

# Update Butan (row 187)
$ws.Range("D187").Value = 225
$ws.Range("E187").Value = 57

# Swap Islas Malvinas (row 215) and Montserrat (row 216):
# the two countries exchange their entire row of data.
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
